# Consolidate Cargo Type & Vehicle Type:
# The VTStFES sheet used to have one row per vehicle type (LDVs, HDVs,
# aircraft, rail, ships, motorbikes). Going forward the model distinguishes
# passenger vs freight cargo, so the existing rows become the "passenger"
# variants and a mirrored set of "freight" rows is appended (initially
# identical, via formulas referencing the passenger row above).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("VTStFES")

# Relabel the existing vehicle-type rows as the "passenger" variants.
$ws2.Range("A2").Value = "passenger LDVs"
$ws2.Range("A3").Value = "passenger HDVs"
$ws2.Range("A4").Value = "passenger aircraft"
$ws2.Range("A5").Value = "passenger rail"
$ws2.Range("A6").Value = "passenger ships"
$ws2.Range("A7").Value = "passenger motorbikes"

# Append the mirrored "freight" rows, each referencing the corresponding
# passenger row's values via formula.
$ws2.Range("A8").Value  = "freight LDVs"
$ws2.Range("A9").Value  = "freight HDVs"
$ws2.Range("A10").Value = "freight aircraft"
$ws2.Range("A11").Value = "freight rail"
$ws2.Range("A12").Value = "freight ships"
$ws2.Range("A13").Value = "freight motorbikes"

$ws2.Range("B8:H8").Formula  = "=B2"
$ws2.Range("B9:H13").Formula = "=B3"

# Widen column A so the longer "passenger"/"freight" labels fit.
$ws2.Columns.Item(1).ColumnWidth = 21.95
